$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D hold numeric-looking text (e.g. "217.26", "1.648.25") that must
# stay text, matching the source inlineStr cells. Force text format, assign,
# then restore the Normal style so no stray style index is introduced.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.148.05"
$ws.Range("E2").Value = "  -0.64%  "
Set-TextValue "D3" "1.643.69"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "217.26"
$ws.Range("E5").Value = "  -1.29%  "
Set-TextValue "D6" "0.509"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.39%  "
Set-TextValue "D10" "19.95"
$ws.Range("E10").Value = "  +0.01%  "
Set-TextValue "D11" "0.0844"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  -1.00%  "
Set-TextValue "D13" "1.648.25"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -2.43%  "
Set-TextValue "D15" "0.539"
$ws.Range("E15").Value = "  +0.77%  "
Set-TextValue "D16" "67.37"
$ws.Range("E16").Value = "  -0.08%  "
Set-TextValue "D17" "27.127.67"
$ws.Range("E17").Value = "  -0.68%  "
Set-TextValue "D18" "0.0₃0739"
$ws.Range("E18").Value = "  +0.34%  "
Set-TextValue "D19" "218.21"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("E20").Value = "  -0.19%  "
Set-TextValue "D21" "6.83"
$ws.Range("E21").Value = "  +1.07%  "
Set-TextValue "D22" "4.44"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  -0.08%  "
Set-TextValue "D24" "9.18"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  -1.09%  "
Set-TextValue "D29" "15.73"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -2.12%  "
Set-TextValue "D31" "1.18"
$ws.Range("E31").Value = "  -1.50%  "
Set-TextValue "D32" "3.37"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +1.21%  "
Set-TextValue "D35" "1.267.25"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  +0.42%  "
Set-TextValue "D39" "0.839"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("E42").Value = "  +4.34%  "
Set-TextValue "D43" "5.40"
$ws.Range("E43").Value = "  +0.01%  "
Set-TextValue "D44" "1.785.59"
$ws.Range("E44").Value = "  -1.13%  "
Set-TextValue "D45" "62.39"
$ws.Range("E45").Value = "  +0.76%  "
Set-TextValue "D46" "91.70"
$ws.Range("E46").Value = "  -0.84%  "
Set-TextValue "D47" "1.60"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E48").Value = "  +17.48%  "
$ws.Range("E49").Value = "  -0.98%  "
Set-TextValue "D50" "7.67"
$ws.Range("E50").Value = "  -0.20%  "
Set-TextValue "D51" "0.0973"
$ws.Range("E51").Value = "  -1.16%  "
